# Refresh cryptocurrency symbol data for columns D (Price) and E (Volume 1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must remain text, matching the
# original inline-string cells, so force a text number format before writing).
$updates = @{
    'D2' = '310.52'
    'E2' = '1.42%'
    'E3' = '-1.94%'
    'D4' = '5.117'
    'E4' = '1.15%'
    'D5' = '0.08201'
    'E5' = '3.35%'
    'D6' = '1.993'
    'E6' = '-10.61%'
    'D7' = '7.969'
    'E7' = '-0.48%'
    'D8' = '2.901'
    'E8' = '10.21%'
    'D9' = '0.9276'
    'E9' = '-0.08%'
    'D10' = '0.1075'
    'E10' = '9.55%'
    'D11' = '0.1932'
    'E11' = '3.00%'
    'D12' = '0.09665'
    'E12' = '3.95%'
    'D13' = '0.03617'
    'E13' = '-2.44%'
    'D14' = '0.09908'
    'D15' = '0.001438'
    'E15' = '0.52%'
    'D16' = '0.005680'
    'E16' = '0.76%'
    'D17' = '3.477'
    'E17' = '0.46%'
    'D18' = '4.129'
    'E18' = '-0.24%'
    'D19' = '0.3419'
    'E19' = '1.42%'
    'D20' = '0.1302'
    'E20' = '-1.23%'
    'D21' = '5.092'
    'E21' = '0.05%'
    'D22' = '0.2192'
    'E22' = '-2.44%'
    'E23' = '-0.31%'
    'D24' = '0.001226'
    'E24' = '-0.90%'
    'E25' = '0.00%'
    'D26' = '0.0001251'
    'E26' = '-3.70%'
    'D27' = '0.0004451'
    'E27' = '-6.07%'
    'D39' = '0.01976'
    'E39' = '3.00%'
    'D40' = '0.04904'
    'E40' = '-0.30%'
    'D41' = '0.007844'
    'E41' = '0.29%'
    'D42' = '0.009687'
    'E42' = '24.24%'
    'E43' = '-1.00%'
    'D44' = '0.002116'
    'E44' = '-4.60%'
    'D45' = '0.01155'
    'E45' = '1.12%'
    'D46' = '0.00006495'
    'E46' = '3.39%'
    'D47' = '0.00000000751'
    'E47' = '0.13%'
    'D48' = '64.38'
    'E48' = '24.38%'
    'D49' = '0.001301'
    'E49' = '-27.69%'
    'D50' = '0.00002102'
    'E50' = '0.13%'
    'D51' = '0.0002001'
    'E51' = '0.13%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
